# Update the weekly fruit/vegetable data: the rows' Fecha (D), Volumen (M),
# Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P) and
# Precio $/Kg (S) values are re-shuffled across rows 2-12 (a re-sort by date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
$rows = @{
    2  = @(44482, 240, 10000, 11000, 10500, 5250)
    3  = @(44818, 200, 11000, 12000, 11500, 5750)
    4  = @(44497, 500,  9000, 10000,  9500, 4750)
    5  = @(44461, 200, 11000, 12000, 11500, 5750)
    6  = @(44475, 240, 11000, 12000, 11500, 5750)
    7  = @(44489, 160,  9500, 10000,  9750, 4875)
    8  = @(44455, 200, 12000, 13000, 12500, 6250)
    9  = @(44819, 240, 11000, 12000, 11500, 5750)
    10 = @(44490, 400,  9500, 10000,  9750, 4875)
    11 = @(44454, 160, 12000, 13000, 12500, 6250)
    12 = @(44517, 400,  5500,  6000,  5750, 2875)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("M$r").Value = $vals[1]
    $ws.Range("N$r").Value = $vals[2]
    $ws.Range("O$r").Value = $vals[3]
    $ws.Range("P$r").Value = $vals[4]
    $ws.Range("S$r").Value = $vals[5]
}
